$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped
# by one day (45171 -> 45172) for every data row (rows 2 through 338).
$range = $ws.Range("C2:C338")
$range.Value = 45172
